# Generate Report for Handback
# Populates the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on the per-locale sheets, flips the
# Status column from "Ready for handoff" to "Handed back: in sync with
# en-US", and widens the columns that now hold the longer handback
# file-name / status strings.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdTarget95 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5229cce256e18f4f9b1c449e282a47093965d00d/e2e/95f1a85d-2ccc-4db5-be14-dc02a2ef1aeb.md"
$mdTargetA1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5229cce256e18f4f9b1c449e282a47093965d00d/e2e/a124d3ff-d0c1-4719-9b83-a01ee0e6541f.md"

$md95 = "95f1a85d-2ccc-4db5-be14-dc02a2ef1aeb.md"
$mdA1 = "a124d3ff-d0c1-4719-9b83-a01ee0e6541f.md"

# --- Overview sheet: refresh the rolled-up status text for both rows ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $statusNew
$ov.Range("F2").Value = $statusNew
$ov.Range("E3").Value = $statusNew
$ov.Range("F3").Value = $statusNew

$ov.Columns.Item(5).ColumnWidth = 29.17
$ov.Columns.Item(6).ColumnWidth = 29.17

# --- Per-locale sheets: zh-cn handed back at 08:54:24, de-de at 08:54:41 ---
$locales = @(
    @{ Name = "zh-cn"; HandbackTime = "2016-10-26 08:54:24" },
    @{ Name = "de-de"; HandbackTime = "2016-10-26 08:54:41" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Name)

    # Status column (C)
    $ws.Range("C2").Value = $statusNew
    $ws.Range("C3").Value = $statusNew

    # Row 2 -> 95f1a85d... file
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdTarget95, "", "", $md95)
    $ws.Range("J2").Value = "95f1a85d-2ccc-4db5-be14-dc02a2ef1aeb.f631695a710692928112b06ffcef2ac5e3f4ba4e." + $locale.Name + ".xlf"
    $ws.Range("K2").Value = $locale.HandbackTime

    # Row 3 -> a124d3ff... file
    $ws.Hyperlinks.Add($ws.Range("I3"), $mdTargetA1, "", "", $mdA1)
    $ws.Range("J3").Value = "a124d3ff-d0c1-4719-9b83-a01ee0e6541f.190d08bc03a9e7ae81afaf570808a64ee9713911." + $locale.Name + ".xlf"
    $ws.Range("K3").Value = $locale.HandbackTime

    # Widen columns to fit the newly-populated long strings
    $ws.Columns.Item(3).ColumnWidth = 29.17
    $ws.Columns.Item(9).ColumnWidth = 39.17
    $ws.Columns.Item(10).ColumnWidth = 39.17
}
